$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6534
$ws.Range("I3").Value = 6813
$ws.Range("C4").Value = 1817
$ws.Range("I4").Value = 1561
$ws.Range("I5").Value = 634
$ws.Range("I6").Value = 7825
$ws.Range("C7").Value = 28360
$ws.Range("I7").Value = 23367

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I5").Value = 69
$ws.Range("I7").Value = 740
$ws.Range("I8").Value = 1396
$ws.Range("I9").Value = 118
$ws.Range("I10").Value = 167
$ws.Range("I11").Value = 356
$ws.Range("I15").Value = 271
$ws.Range("I16").Value = 68
$ws.Range("I18").Value = 177
$ws.Range("I19").Value = 661
$ws.Range("I23").Value = 226
$ws.Range("I29").Value = 1412
$ws.Range("I33").Value = 1046
$ws.Range("I35").Value = 32
$ws.Range("I36").Value = 320
$ws.Range("I37").Value = 736
$ws.Range("I42").Value = 847
$ws.Range("I44").Value = 174
$ws.Range("I47").Value = 168
$ws.Range("I51").Value = 278
$ws.Range("I52").Value = 513
$ws.Range("I55").Value = 266
$ws.Range("C63").Value = 249
$ws.Range("I63").Value = 76
$ws.Range("I65").Value = 536
$ws.Range("I67").Value = 896
$ws.Range("I68").Value = 80
$ws.Range("I78").Value = 315
$ws.Range("I79").Value = 663
$ws.Range("I83").Value = 508
$ws.Range("I84").Value = 207
$ws.Range("I85").Value = 1051
$ws.Range("I87").Value = 53
$ws.Range("I88").Value = 217
$ws.Range("I89").Value = 279
$ws.Range("I90").Value = 302
$ws.Range("I91").Value = 247
$ws.Range("I96").Value = 260
$ws.Range("I98").Value = 164
$ws.Range("I99").Value = 415
$ws.Range("C101").Value = 28360
$ws.Range("I101").Value = 23367

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 298
$ws.Range("I3").Value = 402
$ws.Range("I6").Value = 269
$ws.Range("I7").Value = 1051

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I3").Value = 176
$ws.Range("I6").Value = 146
$ws.Range("I7").Value = 513

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I3").Value = 75
$ws.Range("I7").Value = 356

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 416
$ws.Range("I5").Value = 42
$ws.Range("I7").Value = 1396

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 228
$ws.Range("I6").Value = 198
$ws.Range("I7").Value = 740

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I4").Value = 40
$ws.Range("I7").Value = 279

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I4").Value = 19
$ws.Range("I7").Value = 260

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I6").Value = 221
$ws.Range("I7").Value = 736

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I3").Value = 153
$ws.Range("I7").Value = 415

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 330
$ws.Range("I4").Value = 55
$ws.Range("I7").Value = 896

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 75
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 207

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 179
$ws.Range("I6").Value = 157
$ws.Range("I7").Value = 536

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 184
$ws.Range("I7").Value = 508

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 382
$ws.Range("I7").Value = 1046

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I6").Value = 391
$ws.Range("I7").Value = 1412

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 193
$ws.Range("I6").Value = 210
$ws.Range("I7").Value = 661

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I4").Value = 16
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 174

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 205
$ws.Range("I6").Value = 302
$ws.Range("I7").Value = 847

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 167

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I4").Value = 44
$ws.Range("I7").Value = 315

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 266

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I6").Value = 66
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I2").Value = 77
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 247

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 215
$ws.Range("I6").Value = 196
$ws.Range("I7").Value = 663

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 177

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 92
$ws.Range("I7").Value = 320

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I2").Value = 42
$ws.Range("I7").Value = 168

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 79
$ws.Range("I7").Value = 271

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 104
$ws.Range("I7").Value = 164

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("I6").Value = 13
$ws.Range("I7").Value = 32

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I2").Value = 40
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 217

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I3").Value = 18
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 75
$ws.Range("I7").Value = 302

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I2").Value = 58
$ws.Range("I3").Value = 78
$ws.Range("I7").Value = 278

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("I4").Value = 10
$ws.Range("I7").Value = 80

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I6").Value = 46
$ws.Range("I7").Value = 68
